# Datos Evaluación.xlsx -- add a second "Evaluación 2" sheet with the same
# layout/headers as "Evaluación 1", relabel the question headers, and add
# some sample evaluation rows + light formatting (header style, centered
# body cells, column widths).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$headers = @("Fecha","Gerencia","Cédula","Nombre Líder","Nombre Practicante","Pregunta 1","Pregunta 2","Pregunta 3","Pregunta 4","Pregunta 5","Pregunta 6","Observaciones")

$colWidths = @(10.71, 20.71, 10.71, 16.71, 18.71, 10.71, 10.71, 10.71, 10.71, 10.71, 10.71, 30.71)

function Set-Headers($ws) {
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $cell = $ws.Cells.Item(1, $c + 1)
        $cell.Value = $headers[$c]
        $cell.Font.Bold = $true
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4160
        $cell.Borders.LineStyle = 1
    }
}

function Set-ColumnWidths($ws) {
    for ($c = 0; $c -lt $colWidths.Length; $c++) {
        $ws.Columns.Item($c + 1).ColumnWidth = $colWidths[$c]
    }
}

function Set-TextValue($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Sheet 1: "Evaluación 1" -- relabel headers only, keep existing rows.
# ---------------------------------------------------------------------
Set-Headers $ws1
Set-ColumnWidths $ws1
$ws1.Range("A2:L4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Sheet 2: "Evaluación 2" -- new sheet, same headers + two sample rows.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Evaluación 2"

Set-Headers $ws2
Set-ColumnWidths $ws2

Set-TextValue $ws2 2 1 "2022-11-04"
$ws2.Cells.Item(2, 2).Value = "M&O"
Set-TextValue $ws2 2 3 "13240"
$ws2.Cells.Item(2, 4).Value = "Javier Sarmiento"
$ws2.Cells.Item(2, 5).Value = "Julian Cely"
Set-TextValue $ws2 2 6 "8"
Set-TextValue $ws2 2 7 "7"
Set-TextValue $ws2 2 8 "6"
Set-TextValue $ws2 2 9 "5"
Set-TextValue $ws2 2 10 "4"
Set-TextValue $ws2 2 11 "3"
$ws2.Cells.Item(2, 12).Value = "Que hambre"

Set-TextValue $ws2 3 1 "2022-12-14"
$ws2.Cells.Item(3, 2).Value = "Ejemploooooo oooooooo ooooooooooooooo"
Set-TextValue $ws2 3 3 "10000001231"
$ws2.Cells.Item(3, 4).Value = "Ejemploooooooooo oooooooooo ooooooooo"
$ws2.Cells.Item(3, 5).Value = "Ejemplooooooo ooo oooooo ooooooooooooo"
Set-TextValue $ws2 3 6 "6"
Set-TextValue $ws2 3 7 "7"
Set-TextValue $ws2 3 8 "6"
Set-TextValue $ws2 3 9 "7"
Set-TextValue $ws2 3 10 "7"
Set-TextValue $ws2 3 11 "7"
$ws2.Cells.Item(3, 12).Value = "Aquí va una observación"

$ws2.Range("A2:L3").HorizontalAlignment = -4108

$ws1.Select()
